# Update LR-pair NATMI output values per revised analysis (Dr Hou advice)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 23.496322
$ws.Range("H2").Value = 70.488966
$ws.Range("I2").Value = 0.1321353991144917
$ws.Range("J2").Value = 0.1321353991144917
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 21.08181366666667
$ws.Range("N2").Value = 63.245441
$ws.Range("O2").Value = 0.0571606014598545
$ws.Range("P2").Value = 0.0571606014598545
$ws.Range("Q2").Value = 495.3450822560007
$ws.Range("R2").Value = 4458.105740304006
$ws.Range("S2").Value = 0.007552938887522269
$ws.Range("T2").Value = 0.007552938887522269

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 23.496322
$ws.Range("H3").Value = 70.488966
$ws.Range("I3").Value = 0.1321353991144917
$ws.Range("J3").Value = 0.1321353991144917
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 301.6001486666667
$ws.Range("N3").Value = 904.800446
$ws.Range("O3").Value = 0.8177496571571792
$ws.Range("P3").Value = 0.8177496571571792
$ws.Range("Q3").Value = 7086.494208319871
$ws.Range("R3").Value = 63778.44787487884
$ws.Range("S3").Value = 0.1080536773242026
$ws.Range("T3").Value = 0.1080536773242026

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 23.496322
$ws.Range("H4").Value = 70.488966
$ws.Range("I4").Value = 0.1321353991144917
$ws.Range("J4").Value = 0.1321353991144917
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 46.13524966666667
$ws.Range("N4").Value = 138.405749
$ws.Range("O4").Value = 0.1250897413829664
$ws.Range("P4").Value = 0.1250897413829664
$ws.Range("Q4").Value = 1084.008681718393
$ws.Range("R4").Value = 9756.078135465536
$ws.Range("S4").Value = 0.0165287829027668
$ws.Range("T4").Value = 0.0165287829027668

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 124.9120333333333
$ws.Range("H5").Value = 374.7361
$ws.Range("I5").Value = 0.7024631931202969
$ws.Range("J5").Value = 0.7024631931202969
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 21.08181366666667
$ws.Range("N5").Value = 63.245441
$ws.Range("O5").Value = 0.0571606014598545
$ws.Range("P5").Value = 0.0571606014598545
$ws.Range("Q5").Value = 2633.372211457789
$ws.Range("R5").Value = 23700.3499031201
$ws.Range("S5").Value = 0.0401532186221661
$ws.Range("T5").Value = 0.0401532186221661

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 124.9120333333333
$ws.Range("H6").Value = 374.7361
$ws.Range("I6").Value = 0.7024631931202969
$ws.Range("J6").Value = 0.7024631931202969
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 301.6001486666667
$ws.Range("N6").Value = 904.800446
$ws.Range("O6").Value = 0.8177496571571792
$ws.Range("P6").Value = 0.8177496571571792
$ws.Range("Q6").Value = 37673.48782358896
$ws.Range("R6").Value = 339061.3904123005
$ws.Range("S6").Value = 0.5744390353396601
$ws.Range("T6").Value = 0.5744390353396601

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 124.9120333333333
$ws.Range("H7").Value = 374.7361
$ws.Range("I7").Value = 0.7024631931202969
$ws.Range("J7").Value = 0.7024631931202969
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 46.13524966666667
$ws.Range("N7").Value = 138.405749
$ws.Range("O7").Value = 0.1250897413829664
$ws.Range("P7").Value = 0.1250897413829664
$ws.Range("Q7").Value = 5762.847844204322
$ws.Range("R7").Value = 51865.6305978389
$ws.Range("S7").Value = 0.08787093915847069
$ws.Range("T7").Value = 0.08787093915847069

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 29.411685
$ws.Range("H8").Value = 88.235055
$ws.Range("I8").Value = 0.1654014077652114
$ws.Range("J8").Value = 0.1654014077652114
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 21.08181366666667
$ws.Range("N8").Value = 63.245441
$ws.Range("O8").Value = 0.0571606014598545
$ws.Range("P8").Value = 0.0571606014598545
$ws.Range("Q8").Value = 620.051662792695
$ws.Range("R8").Value = 5580.464965134255
$ws.Range("S8").Value = 0.009454443950166131
$ws.Range("T8").Value = 0.009454443950166131

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 29.411685
$ws.Range("H9").Value = 88.235055
$ws.Range("I9").Value = 0.1654014077652114
$ws.Range("J9").Value = 0.1654014077652114
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 301.6001486666667
$ws.Range("N9").Value = 904.800446
$ws.Range("O9").Value = 0.8177496571571792
$ws.Range("P9").Value = 0.8177496571571792
$ws.Range("Q9").Value = 8870.56856853717
$ws.Range("R9").Value = 79835.11711683453
$ws.Range("S9").Value = 0.1352569444933164
$ws.Range("T9").Value = 0.1352569444933164

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 29.411685
$ws.Range("H10").Value = 88.235055
$ws.Range("I10").Value = 0.1654014077652114
$ws.Range("J10").Value = 0.1654014077652114
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 46.13524966666667
$ws.Range("N10").Value = 138.405749
$ws.Range("O10").Value = 0.1250897413829664
$ws.Range("P10").Value = 0.1250897413829664
$ws.Range("Q10").Value = 1356.915430592355
$ws.Range("R10").Value = 12212.2388753312
$ws.Range("S10").Value = 0.02069001932172886
$ws.Range("T10").Value = 0.02069001932172886

